{"js": "const body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\n// Map of literal prefix -> placeholder field name. Each matching\n// paragraph's \"Label: value\" text is rewritten to \"Label: {Field}\".\nconst replacements = [\n  { prefix: \"Application ID: \", field: \"ApplicationID\" },\n  { prefix: \"Application Name: \", field: \"ApplicationName\" },\n  { prefix: \"Owner: \", field: \"Owner\" },\n  { prefix: \"Technology: \", field: \"Technology\" },\n  { prefix: \"Hosting: \", field: \"Hosting\" },\n  { prefix: \"Database: \", field: \"Database\" },\n  { prefix: \"Integrations: \", field: \"Integrations\" },\n  { prefix: \"Description: \", field: \"Description\" },\n  { prefix: \"Risks: \", field: \"Risks\" }\n];\n\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const para = paragraphs.items[i];\n  const text = para.text;\n  for (const { prefix, field } of replacements) {\n    if (text.indexOf(prefix) === 0) {\n      const newText = prefix + \"{\" + field + \"}\";\n      if (newText !== text) {\n        // Clear the paragraph's content first, then insert the fresh text.\n        // (Editing the run's text in place keeps stale run-level XML quirks\n        // such as a lingering `xml:space=\"preserve\"`; clearing + inserting\n        // mirrors how Word rewrites the run from scratch.)\n        para.clear();\n        para.insertText(newText, \"Start\");\n      }\n      break;\n    }\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace each \"Label: value\" line's value with a \"{FieldName}\" placeholder,\n# matching the literal old text exactly so only the intended paragraphs change.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @{ Old = \"Application ID: APP002\";       New = \"Application ID: {ApplicationID}\" },\n    @{ Old = \"Application Name: HRPortal\";    New = \"Application Name: {ApplicationName}\" },\n    @{ Old = \"Owner: Mary\";                   New = \"Owner: {Owner}\" },\n    @{ Old = \"Technology: Java/Spring\";       New = \"Technology: {Technology}\" },\n    @{ Old = \"Hosting: Azure\";                New = \"Hosting: {Hosting}\" },\n    @{ Old = \"Database: MySQL\";               New = \"Database: {Database}\" },\n    @{ Old = \"Integrations: Workday\";         New = \"Integrations: {Integrations}\" },\n    @{ Old = \"Description: \";                 New = \"Description: {Description}\" },\n    @{ Old = \"Risks: \";                       New = \"Risks: {Risks}\" }\n)\n\nforeach ($pair in $replacements) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $pair.Old\n    $find.Replacement.Text = $pair.New\n    $find.Execute([ref]$pair.Old, $false, $false, $false, $false, $false, $true, 1, $false, [ref]$pair.New, 2) | Out-Null\n}\n"}
